$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.647.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "'1.642.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'214.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "'0.0625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'1.872.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "'1.650.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "'4.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "'64.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'26.686.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "'0.0₃0742"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "'215.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'1.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'4.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").Value = "'6.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'9.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").Value = "'2.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.39%  "
$ws.Range("D25").Value = "'145.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'7.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").Value = "'15.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "'0.0513"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("D34").Value = "'1.279.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("D38").Value = "'0.530"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("D39").Value = "'0.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'0.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D42").Value = "'2.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  +1.27%  "
$ws.Range("D44").Value = "'1.782.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "'91.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'59.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.30%  "
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'7.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").Value = "'0.0963"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("D51").Value = "'0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
